# Daily attendance processing - 2026-01-02 19:27:28
# Swap the order of names in the "Recorded By" column (G) wherever the
# value is exactly "dnasr281@gmail.com, System", turning it into
# "System, dnasr281@gmail.com".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Session Analysis Results")

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

$oldValue = "dnasr281@gmail.com, System"
$newValue = "System, dnasr281@gmail.com"

for ($row = 1; $row -le $lastRow; $row++) {
    $cell = $ws.Cells.Item($row, 7)  # Column G
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
